# Swap the data held in rows 5 and 6 of the "Artfynd" sheet.
# (Row 5 becomes the former row 6's record and vice versa.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 gets the values that used to live in row 6 ---
$ws.Cells.Item(5, 1).Value2  = 131067826        # A5  Id
$ws.Cells.Item(5, 2).Value2  = 79243             # B5  Taxonsorteringsordning
$ws.Cells.Item(5, 5).Value2  = 6425              # E5  TaxonId
$ws.Cells.Item(5, 6).Value2  = "Garnlav"         # F5  Artnamn
$ws.Cells.Item(5, 7).Value2  = "Alectoria sarmentosa"  # G5  Vetenskapligt namn
$ws.Cells.Item(5, 8).Value2  = "(Ach.) Ach."     # H5  Auktor
$ws.Cells.Item(5, 16).Value2 = "Långan Öst, Jmt" # P5  Lokalnamn
$ws.Cells.Item(5, 17).Value2 = 465891            # Q5  Ost
$ws.Cells.Item(5, 18).Value2 = 7046290           # R5  Nord
$ws.Cells.Item(5, 34).Value2 = "Granskog"        # AH5 Biotop
$ws.Cells.Item(5, 49).Value2 = "Kristian Zackrisson"  # AW5 Rapportör
$ws.Cells.Item(5, 50).Value2 = "Kristian Zackrisson"  # AX5 Observatörer

# Row 5 loses the values that were specific to the old row 5 record
$ws.Cells.Item(5, 12).ClearContents()  # L5  Kön   (old row 5 had a blank tag here)
$ws.Cells.Item(5, 13).ClearContents()  # M5  Aktivitet
$ws.Cells.Item(5, 26).ClearContents()  # Z5  Starttid
$ws.Cells.Item(5, 28).ClearContents()  # AB5 Sluttid
$ws.Cells.Item(5, 29).ClearContents()  # AC5 Publik kommentar

# Row 5 picks up the (empty) placeholder tags that row 6 used to carry.
# Touching a no-op formatting property is enough to keep the blank cell
# without putting any value or visible formatting in it.
$ws.Cells.Item(5, 10).Font.Bold = $false  # J5  Enhet
$ws.Cells.Item(5, 32).Font.Bold = $false  # AF5 Bestämningsmetod

# --- Row 6 gets the values that used to live in row 5 ---
$ws.Cells.Item(6, 1).Value2  = 131067473        # A6  Id
$ws.Cells.Item(6, 2).Value2  = 57884             # B6  Taxonsorteringsordning
$ws.Cells.Item(6, 5).Value2  = 100109            # E6  TaxonId
$ws.Cells.Item(6, 6).Value2  = "Tretåig hackspett"     # F6  Artnamn
$ws.Cells.Item(6, 7).Value2  = "Picoides tridactylus"  # G6  Vetenskapligt namn
$ws.Cells.Item(6, 8).Value2  = "(Linnaeus, 1758)"      # H6  Auktor
$ws.Cells.Item(6, 13).Value2 = "färska spår"     # M6  Aktivitet
$ws.Cells.Item(6, 16).Value2 = "Åbogen, Jmt"     # P6  Lokalnamn
$ws.Cells.Item(6, 17).Value2 = 465809            # Q6  Ost
$ws.Cells.Item(6, 18).Value2 = 7046259           # R6  Nord
$ws.Cells.Item(6, 26).Value2 = "15:46"           # Z6  Starttid
$ws.Cells.Item(6, 28).Value2 = "15:46"           # AB6 Sluttid
$ws.Cells.Item(6, 29).Value2 = "Färska ringhack" # AC6 Publik kommentar
$ws.Cells.Item(6, 49).Value2 = "Elin Albrechtsson"     # AW6 Rapportör
$ws.Cells.Item(6, 50).Value2 = "Elin Albrechtsson"     # AX6 Observatörer

# Row 6 loses the values that were specific to the old row 6 record
$ws.Cells.Item(6, 10).ClearContents()  # J6  Enhet
$ws.Cells.Item(6, 32).ClearContents()  # AF6 Bestämningsmetod
$ws.Cells.Item(6, 34).ClearContents()  # AH6 Biotop

# Row 6 picks up the (empty) placeholder tag that row 5 used to carry.
$ws.Cells.Item(6, 12).Font.Bold = $false  # L6  Kön
